$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Cells.Item(98, 8).Value = 4108.154
$ws.Cells.Item(98, 9).Value = 1662.5
$ws.Cells.Item(98, 10).Value = 8021.2
$ws.Cells.Item(98, 11).Value = 1662.5
$ws.Cells.Item(98, 12).Value = 8021.2
$ws.Cells.Item(98, 13).Value = -164.5
$ws.Cells.Item(98, 14).Value = -11017.2

# Row 122
$ws.Cells.Item(122, 8).Value = 4108.154
$ws.Cells.Item(122, 9).Value = 1662.5
$ws.Cells.Item(122, 10).Value = 8021.2
$ws.Cells.Item(122, 11).Value = 4987.5
$ws.Cells.Item(122, 12).Value = 24063.6
$ws.Cells.Item(122, 13).Value = -2537.5
$ws.Cells.Item(122, 14).Value = -28963.6

# Row 129
$ws.Cells.Item(129, 8).Value = 915.6087
$ws.Cells.Item(129, 9).Value = 482.22223
$ws.Cells.Item(129, 10).Value = 1194.2142
$ws.Cells.Item(129, 11).Value = 1446.66669
$ws.Cells.Item(129, 12).Value = 3582.6426
$ws.Cells.Item(129, 13).Value = 3553.33331
$ws.Cells.Item(129, 14).Value = -13582.6426

# Row 132
$ws.Cells.Item(132, 8).Value = 2835.077
$ws.Cells.Item(132, 9).Value = 869.95654
$ws.Cells.Item(132, 10).Value = 17901
$ws.Cells.Item(132, 11).Value = 2609.86962
$ws.Cells.Item(132, 12).Value = 53703
$ws.Cells.Item(132, 13).Value = -79.86961999999994
$ws.Cells.Item(132, 14).Value = -58763

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 15417.07
$ws.Cells.Item(32, 9).Value = 7090.4854
$ws.Cells.Item(32, 10).Value = 33111.062
$ws.Cells.Item(32, 11).Value = 7090.4854
$ws.Cells.Item(32, 12).Value = 33111.062
$ws.Cells.Item(32, 13).Value = -6803.4854
$ws.Cells.Item(32, 14).Value = -33685.062

# Row 110
$ws.Cells.Item(110, 8).Value = 1207.6
$ws.Cells.Item(110, 9).Value = 1073.625
$ws.Cells.Item(110, 10).Value = 1360.7142
$ws.Cells.Item(110, 11).Value = 1073.625
$ws.Cells.Item(110, 12).Value = 1360.7142
$ws.Cells.Item(110, 13).Value = 971.375
$ws.Cells.Item(110, 14).Value = -5450.7142

# Row 114
$ws.Cells.Item(114, 8).Value = 28778.6
$ws.Cells.Item(114, 10).Value = 28778.6
$ws.Cells.Item(114, 12).Value = 28778.6
$ws.Cells.Item(114, 14).Value = -37456.6

# Row 122
$ws.Cells.Item(122, 8).Value = 1763.6072
$ws.Cells.Item(122, 9).Value = 1463.1666
$ws.Cells.Item(122, 10).Value = 2304.4
$ws.Cells.Item(122, 11).Value = 4389.4998
$ws.Cells.Item(122, 12).Value = 6913.200000000001
$ws.Cells.Item(122, 13).Value = -1939.4998
$ws.Cells.Item(122, 14).Value = -11813.2

# Row 132
$ws.Cells.Item(132, 8).Value = 2230.2856
$ws.Cells.Item(132, 9).Value = 2376.25
$ws.Cells.Item(132, 10).Value = 2140.4614
$ws.Cells.Item(132, 11).Value = 7128.75
$ws.Cells.Item(132, 12).Value = 6421.3842
$ws.Cells.Item(132, 13).Value = -4598.75
$ws.Cells.Item(132, 14).Value = -11481.3842

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 1589.9546
$ws.Cells.Item(86, 10).Value = 1391.5
$ws.Cells.Item(86, 12).Value = 1391.5
$ws.Cells.Item(86, 14).Value = -3637.5

# Row 89
$ws.Cells.Item(89, 8).Value = 1589.9546
$ws.Cells.Item(89, 10).Value = 1391.5
$ws.Cells.Item(89, 12).Value = 6957.5
$ws.Cells.Item(89, 14).Value = -18189.5

# Row 105
$ws.Cells.Item(105, 8).Value = 2406.8293
$ws.Cells.Item(105, 9).Value = 2356.2163
$ws.Cells.Item(105, 10).Value = 2875
$ws.Cells.Item(105, 11).Value = 2356.2163
$ws.Cells.Item(105, 12).Value = 2875
$ws.Cells.Item(105, 13).Value = -609.2163
$ws.Cells.Item(105, 14).Value = -6369

# Row 134
$ws.Cells.Item(134, 8).Value = 793661.5600000001
$ws.Cells.Item(134, 9).Value = 1029329.1
$ws.Cells.Item(134, 10).Value = 8102.8335
$ws.Cells.Item(134, 11).Value = 3087987.3
$ws.Cells.Item(134, 12).Value = 24308.5005
$ws.Cells.Item(134, 13).Value = -3085452.3
$ws.Cells.Item(134, 14).Value = -29378.5005

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2623.2068
$ws.Cells.Item(31, 9).Value = 1210.3823
$ws.Cells.Item(31, 10).Value = 4624.7085
$ws.Cells.Item(31, 11).Value = 1210.3823
$ws.Cells.Item(31, 12).Value = 4624.7085
$ws.Cells.Item(31, 13).Value = -915.3823
$ws.Cells.Item(31, 14).Value = -5214.7085

# Row 34
$ws.Cells.Item(34, 8).Value = 2623.2068
$ws.Cells.Item(34, 9).Value = 1210.3823
$ws.Cells.Item(34, 10).Value = 4624.7085
$ws.Cells.Item(34, 11).Value = 1210.3823
$ws.Cells.Item(34, 12).Value = 4624.7085
$ws.Cells.Item(34, 13).Value = -1008.3823
$ws.Cells.Item(34, 14).Value = -5028.7085

# Row 99
$ws.Cells.Item(99, 8).Value = 2671.5
$ws.Cells.Item(99, 9).Value = 2671.5
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 2671.5
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -1173.5
$ws.Cells.Item(99, 14).ClearContents()

# Row 126
$ws.Cells.Item(126, 8).Value = 2671.5
$ws.Cells.Item(126, 9).Value = 2671.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8014.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -5544.5
$ws.Cells.Item(126, 14).ClearContents()

# Row 135
$ws.Cells.Item(135, 8).Value = 37285.715
$ws.Cells.Item(135, 10).Value = 38700
$ws.Cells.Item(135, 12).Value = 38700
$ws.Cells.Item(135, 14).Value = -48840

# Row 141
$ws.Cells.Item(141, 8).Value = 40948.7
$ws.Cells.Item(141, 9).Value = 19648
$ws.Cells.Item(141, 10).Value = 43315.445
$ws.Cells.Item(141, 11).Value = 19648
$ws.Cells.Item(141, 12).Value = 43315.445
$ws.Cells.Item(141, 13).Value = -14468
$ws.Cells.Item(141, 14).Value = -53675.445

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Cells.Item(107, 8).Value = 5555970
$ws.Cells.Item(107, 9).Value = 497.25
$ws.Cells.Item(107, 10).Value = 16666915
$ws.Cells.Item(107, 11).Value = 1491.75
$ws.Cells.Item(107, 12).Value = 50000745
$ws.Cells.Item(107, 13).Value = 428.25
$ws.Cells.Item(107, 14).Value = -50004585

# Row 113
$ws.Cells.Item(113, 8).Value = 2020761.5
$ws.Cells.Item(113, 9).Value = 2755379.5
$ws.Cells.Item(113, 10).Value = 562
$ws.Cells.Item(113, 11).Value = 8266138.5
$ws.Cells.Item(113, 12).Value = 1686
$ws.Cells.Item(113, 13).Value = -8263968.5
$ws.Cells.Item(113, 14).Value = -6026

# Row 132
$ws.Cells.Item(132, 8).Value = 626547.7
$ws.Cells.Item(132, 9).Value = 1621.4445
$ws.Cells.Item(132, 10).Value = 1924471.5
$ws.Cells.Item(132, 11).Value = 14593.0005
$ws.Cells.Item(132, 12).Value = 17320243.5
$ws.Cells.Item(132, 13).Value = -12063.0005
$ws.Cells.Item(132, 14).Value = -17325303.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 2180
$ws.Cells.Item(80, 9).Value = 2328.5715
$ws.Cells.Item(80, 10).Value = 1833.3334
$ws.Cells.Item(80, 11).Value = 2328.5715
$ws.Cells.Item(80, 12).Value = 1833.3334
$ws.Cells.Item(80, 13).Value = -1330.5715
$ws.Cells.Item(80, 14).Value = -3829.3334

# Row 83
$ws.Cells.Item(83, 8).Value = 2180
$ws.Cells.Item(83, 9).Value = 2328.5715
$ws.Cells.Item(83, 10).Value = 1833.3334
$ws.Cells.Item(83, 11).Value = 11642.8575
$ws.Cells.Item(83, 12).Value = 9166.666999999999
$ws.Cells.Item(83, 13).Value = -6650.8575
$ws.Cells.Item(83, 14).Value = -19150.667

# Row 102
$ws.Cells.Item(102, 8).Value = 5556977
$ws.Cells.Item(102, 9).Value = 10102058
$ws.Cells.Item(102, 11).Value = 10102058
$ws.Cells.Item(102, 13).Value = -10100436

# Row 103
$ws.Cells.Item(103, 8).Value = 24650
$ws.Cells.Item(103, 10).Value = 24650
$ws.Cells.Item(103, 12).Value = 24650
$ws.Cells.Item(103, 14).Value = -26994

# Row 123
$ws.Cells.Item(123, 8).Value = 15326
$ws.Cells.Item(123, 10).Value = 15326
$ws.Cells.Item(123, 12).Value = 15326
$ws.Cells.Item(123, 14).Value = -20226

# Row 132
$ws.Cells.Item(132, 8).Value = 1671064.8
$ws.Cells.Item(132, 9).Value = 3481.3635
$ws.Cells.Item(132, 10).Value = 3709222.2
$ws.Cells.Item(132, 11).Value = 10444.0905
$ws.Cells.Item(132, 12).Value = 11127666.6
$ws.Cells.Item(132, 13).Value = -7914.0905
$ws.Cells.Item(132, 14).Value = -11132726.6

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 1600.3334
$ws.Cells.Item(81, 9).Value = 1600.3334
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 3200.6668
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -2139.6668
$ws.Cells.Item(81, 14).ClearContents()

# Row 84
$ws.Cells.Item(84, 8).Value = 1600.3334
$ws.Cells.Item(84, 9).Value = 1600.3334
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 16003.334
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -10699.334
$ws.Cells.Item(84, 14).ClearContents()
